# Generate Report for Handoff
# Update the "Latest Handoff Datetime" column (D) for the file rows that
# were just re-handed-off, on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhRows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $zhRows) {
    $zhcn.Cells.Item($r, 4).Value = "2016-03-08 07:15:50"
}

$deRows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $deRows) {
    $dede.Cells.Item($r, 4).Value = "2016-03-08 07:16:01"
}
